# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Durazno"
# at row 165, shifting the existing rows 165-170 down to 166-171 (so the
# sheet's used range grows from A1:T170 to A1:T171).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 165 (and everything below it) down by one row.
$ws.Rows(165).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A165").Value = 4
$ws.Range("B165").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C165").Value = "Los Lagos"
$ws.Range("D165").Value = 44516
$ws.Range("E165").Value = 10
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100103
$ws.Range("H165").Value = "Frutos de hueso (carozo)"
$ws.Range("I165").Value = 100103004
$ws.Range("J165").Value = "Durazno"
$ws.Range("K165").Value = "Florida King"
$ws.Range("L165").Value = "Tercera"
$ws.Range("M165").Value = 400
$ws.Range("N165").Value = 16000
$ws.Range("O165").Value = 16500
$ws.Range("P165").Value = 16250
$ws.Range("Q165").Value = "$/caja 14 kilos empedrada"
$ws.Range("R165").Value = "Provincia de Limarí"
$ws.Range("S165").Value = 1161
$ws.Range("T165").Value = 14
